# Weekly fruit/vegetable price update: a new weekly record is inserted
# at row 21 (Sandia, "Primera" quality, week of 2022-12-20), pushing all
# subsequent rows down by one. The sheet's dimension grows from R53 to R54.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 21, shifting rows 21:53 down to 22:54
# (and carrying the date-format style on column D along with it).
$ws.Rows.Item(21).Insert()

# Populate the newly inserted row 21 with the new weekly record.
$ws.Cells.Item(21, 1).Value = 1
$ws.Cells.Item(21, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(21, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(21, 4).Value = 44915
$ws.Cells.Item(21, 5).Value = 15
$ws.Cells.Item(21, 6).Value = 100112028
$ws.Cells.Item(21, 7).Value = "Sandia"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 1000
$ws.Cells.Item(21, 11).Value = 440
$ws.Cells.Item(21, 12).Value = 450
$ws.Cells.Item(21, 13).Value = 444
$ws.Cells.Item(21, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(21, 15).Value = "Perú"
$ws.Cells.Item(21, 16).Value = 444
$ws.Cells.Item(21, 17).Value = 1
$ws.Cells.Item(21, 18).Value = "Hortaliza"
